$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "regression" (3rd sheet / sheet3.xml)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

# Flip Run Flag from N -> Y for rows 3 through 20 (row 2 already Y, untouched)
for ($r = 3; $r -le 20; $r++) {
    $ws3.Range("B" + $r).Value = "Y"
}

# Copy cell level formatting for the two brand new rows (21 & 22) from the
# most similar existing rows so the style indices line up with the target.
$ws3.Range("C20").Copy()
$ws3.Range("C21").PasteSpecial(-4122)   # xlPasteFormats
$ws3.Range("D20").Copy()
$ws3.Range("D21").PasteSpecial(-4122)

$ws3.Range("B19").Copy()
$ws3.Range("B22").PasteSpecial(-4122)
$ws3.Range("C19").Copy()
$ws3.Range("C22").PasteSpecial(-4122)
$ws3.Range("D19").Copy()
$ws3.Range("D22").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# New row 21 - CELEBI validate test case
$ws3.Range("A21").Value = "ZestIOT_AV_2307_UI_CELEBI_Validate_LANDING_ONBLOCK_OFFBLOCK_AIRBORNE_timestamps_of_Arrival_and_Departure_aircrafts_Any_Data_source"
$ws3.Range("B21").Value = "Y"
$ws3.Range("C21").Value = "AV 2307 Validate LANDING ONBLOCK OFFBLOCK AIRBORNE timestamps of Arrival and Departure aircrafts Any Data source"
$ws3.Range("D21").Value = "CELEBI-Delhi"

# New row 22 - BSSPL fueling coverage test case
$ws3.Range("A22").Value = "ZestIOT_AV_2405_UI_BSSPL_Validate_Fueling_Coverage"
$ws3.Range("B22").Value = "Y"
$ws3.Range("C22").Value = "AV 2405 Read scheduled flights and check the coverage (Scheduled Flights Vs Flights detected, Scheduled Flights Vs fuel activity formed) - BSSPL"
$ws3.Range("D22").Value = "BSSPL-Delhi"

# Column E ("Execution Type") for the two new rows - added last so the new
# unique string "UI" lands at the expected position in the shared strings table.
$ws3.Range("E21").Value = "UI"
$ws3.Range("E22").Value = "UI"

$ws3.Activate()
$ws3.Range("A22").Select()
$excel.ActiveWindow.ScrollRow = 16

# ---------------------------------------------------------------------------
# Sheet "AppControl" (1st sheet / sheet1.xml)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Swap Run Flag values between "smoke" (row 3) and "regression" (row 4)
$ws1.Range("B3").Value = "N"
$ws1.Range("B4").Value = "Y"

# Update the notification email address
$ws1.Range("B25").Value = "stiyyagura@enhops.com"

$ws1.Activate()
$ws1.Range("B5").Select()
